$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected; unprotect it so the cell values below can
# be written, then re-protect it when we're done.
$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure blurb
# (stored in A11) from 2021-05-05 to 2021-05-06.
$line1 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution."
$line2 = "Model holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."
$ws.Range("A11").Value = $line1 + [char]10 + $line2

# Refresh the Weight (D) and Percent Change (E) columns for the holdings
# table (rows 2-8).
$ws.Range("D2").Value = 0.4977671876148553
$ws.Range("E2").Value = 0.008911091608722055

$ws.Range("D3").Value = 0.2452674095691414
$ws.Range("E3").Value = 0.007171081516171407

$ws.Range("D4").Value = 0.09690151816731629
$ws.Range("E4").Value = -0.0002482621648459693

$ws.Range("D5").Value = 0.1030345208126295
$ws.Range("E5").Value = 0.009357917168535179

$ws.Range("D6").Value = 0.02990511812263905
$ws.Range("E6").Value = 0.009097873970503878

$ws.Range("D7").Value = 0.02712424571341835
$ws.Range("E7").Value = 0.004103755323267455

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0.007517997390008357

# Restore sheet protection to its original (protected) state.
$ws.Protect()
